$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-12 14:34:15"
$wsZhCn.Range("H4").Value = "2016-03-12 14:34:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-12 14:34:18"
$wsDeDe.Range("H4").Value = "2016-03-12 14:34:38"
